$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new data rows (22-30) following the same pattern as the existing
# rows: regcntr_id, usr_id, machine_id, lang_code ("eng"), is_active (TRUE),
# cr_by ("superadmin"), cr_dtimes ("now()").
$newRows = @(
    @(10002, 110021, 10021),
    @(10003, 110022, 10022),
    @(10004, 110023, 10023),
    @(10005, 110024, 10024),
    @(10006, 110025, 10025),
    @(10007, 110026, 10026),
    @(10008, 110027, 10027),
    @(10009, 110028, 10028),
    @(10010, 110029, 10029)
)

$r = 22
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $r = $r + 1
}

# Mirror the selection Excel leaves on the sheet after entering data up to
# row 30 - the full remaining rows below the data get selected.
$ws.Range("A31:XFD1048576").Select()

# Best-effort match of the recorded page setup (orientation = portrait).
$ws.PageSetup.Orientation = 1
